$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3158
$ws.Range("J69").Value = 3158
$ws.Range("L69").Value = 9474
$ws.Range("N69").Value = -11222
$ws.Range("H72").Value = 3158
$ws.Range("J72").Value = 3158
$ws.Range("L72").Value = 28422
$ws.Range("N72").Value = -37158
$ws.Range("H112").Value = 1088.7407
$ws.Range("J112").Value = 1136.2727
$ws.Range("L112").Value = 3408.8181
$ws.Range("N112").Value = -5624.8181
$ws.Range("H113").Value = 1767.1034
$ws.Range("I113").Value = 1112.5
$ws.Range("J113").Value = 2016.4762
$ws.Range("K113").Value = 1112.5
$ws.Range("L113").Value = 2016.4762
$ws.Range("M113").Value = 2141.5
$ws.Range("N113").Value = -8524.476200000001
$ws.Range("H137").Value = 1546.7142
$ws.Range("I137").Value = 1447.2354
$ws.Range("J137").Value = 1700.4546
$ws.Range("K137").Value = 4341.706200000001
$ws.Range("L137").Value = 5101.3638
$ws.Range("M137").Value = -1791.706200000001
$ws.Range("N137").Value = -10201.3638
$ws.Range("H138").Value = 2983.6575
$ws.Range("I138").Value = 2132.2666
$ws.Range("J138").Value = 3203.8447
$ws.Range("K138").Value = 6396.7998
$ws.Range("L138").Value = 9611.534100000001
$ws.Range("M138").Value = -1256.7998
$ws.Range("N138").Value = -19891.5341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1647.76
$ws.Range("I45").Value = 1862.909
$ws.Range("J45").Value = 1478.7142
$ws.Range("K45").Value = 1862.909
$ws.Range("L45").Value = 1478.7142
$ws.Range("M45").Value = -1485.909
$ws.Range("N45").Value = -2232.7142
$ws.Range("H63").Value = 2991.7144
$ws.Range("I63").Value = 2138.4
$ws.Range("J63").Value = 5125
$ws.Range("K63").Value = 2138.4
$ws.Range("L63").Value = 5125
$ws.Range("M63").Value = -1452.4
$ws.Range("N63").Value = -6497
$ws.Range("H66").Value = 2991.7144
$ws.Range("I66").Value = 2138.4
$ws.Range("J66").Value = 5125
$ws.Range("K66").Value = 10692
$ws.Range("L66").Value = 25625
$ws.Range("M66").Value = -7260
$ws.Range("N66").Value = -32489
$ws.Range("H110").Value = 2218.6155
$ws.Range("I110").Value = 2078.7
$ws.Range("J110").Value = 2685
$ws.Range("K110").Value = 2078.7
$ws.Range("L110").Value = 2685
$ws.Range("M110").Value = -33.69999999999982
$ws.Range("N110").Value = -6775

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1002140
$ws.Range("I86").Value = 2501150
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 2501150
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -2500027
$ws.Range("N86").Value = -5046
$ws.Range("H89").Value = 1002140
$ws.Range("I89").Value = 2501150
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 12505750
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -12500134
$ws.Range("N89").Value = -25232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 995.52
$ws.Range("I114").Value = 1333.2142
$ws.Range("J114").Value = 565.7273
$ws.Range("K114").Value = 3999.6426
$ws.Range("L114").Value = 1697.1819
$ws.Range("M114").Value = -745.6425999999997
$ws.Range("N114").Value = -8205.1819
$ws.Range("H130").Value = 1231.125
$ws.Range("I130").Value = 1012.25
$ws.Range("J130").Value = 1450
$ws.Range("K130").Value = 3036.75
$ws.Range("L130").Value = 4350
$ws.Range("M130").Value = 1983.25
$ws.Range("N130").Value = -14390
$ws.Range("H131").Value = 902.59375
$ws.Range("I131").Value = 342.2857
$ws.Range("J131").Value = 1059.48
$ws.Range("K131").Value = 1026.8571
$ws.Range("L131").Value = 3178.44
$ws.Range("M131").Value = 4013.1429
$ws.Range("N131").Value = -13258.44

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1303.5714
$ws.Range("I97").Value = 1324
$ws.Range("J97").Value = 1252.5
$ws.Range("K97").Value = 1324
$ws.Range("L97").Value = 1252.5
$ws.Range("M97").Value = -828
$ws.Range("N97").Value = -2244.5
$ws.Range("H132").Value = 35716530
$ws.Range("I132").Value = 66667932
$ws.Range("J132").Value = 3371.1538
$ws.Range("K132").Value = 200003796
$ws.Range("L132").Value = 10113.4614
$ws.Range("M132").Value = -200001266
$ws.Range("N132").Value = -15173.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2314.7334
$ws.Range("I61").Value = 1931.1
$ws.Range("J61").Value = 3082
$ws.Range("K61").Value = 1931.1
$ws.Range("L61").Value = 3082
$ws.Range("M61").Value = -1729.1
$ws.Range("N61").Value = -3486
$ws.Range("H68").Value = 17809722
$ws.Range("I68").Value = 37594716
$ws.Range("J68").Value = 3227.2
$ws.Range("K68").Value = 37594716
$ws.Range("L68").Value = 3227.2
$ws.Range("M68").Value = -37593967
$ws.Range("N68").Value = -4725.2
$ws.Range("H71").Value = 17809722
$ws.Range("I71").Value = 37594716
$ws.Range("J71").Value = 3227.2
$ws.Range("K71").Value = 187973580
$ws.Range("L71").Value = 16136
$ws.Range("M71").Value = -187969836
$ws.Range("N71").Value = -23624
$ws.Range("H82").Value = 1955.3636
$ws.Range("I82").Value = 1321.8
$ws.Range("J82").Value = 2483.3333
$ws.Range("K82").Value = 1321.8
$ws.Range("L82").Value = 2483.3333
$ws.Range("M82").Value = -960.8
$ws.Range("N82").Value = -3205.3333
$ws.Range("H85").Value = 1955.3636
$ws.Range("I85").Value = 1321.8
$ws.Range("J85").Value = 2483.3333
$ws.Range("K85").Value = 1321.8
$ws.Range("L85").Value = 2483.3333
$ws.Range("M85").Value = -73.79999999999995
$ws.Range("N85").Value = -4979.3333
$ws.Range("H111").Value = 35001
$ws.Range("J111").Value = 35001
$ws.Range("L111").Value = 35001
$ws.Range("N111").Value = -43181
$ws.Range("H113").Value = 2314.7334
$ws.Range("I113").Value = 1931.1
$ws.Range("J113").Value = 3082
$ws.Range("K113").Value = 1931.1
$ws.Range("L113").Value = 3082
$ws.Range("M113").Value = 238.9000000000001
$ws.Range("N113").Value = -7422
$ws.Range("H122").Value = 4583.5
$ws.Range("I122").Value = 4157
$ws.Range("J122").Value = 6076.25
$ws.Range("K122").Value = 12471
$ws.Range("L122").Value = 18228.75
$ws.Range("M122").Value = -10021
$ws.Range("N122").Value = -23128.75
$ws.Range("H136").Value = 9263993
$ws.Range("I136").Value = 22231038
$ws.Range("J136").Value = 1817.619
$ws.Range("K136").Value = 66693114
$ws.Range("L136").Value = 5452.857
$ws.Range("M136").Value = -66690564
$ws.Range("N136").Value = -10552.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 22639.814
$ws.Range("I100").Value = 42820.355
$ws.Range("K100").Value = 85640.71000000001
$ws.Range("M100").Value = -85099.71000000001
$ws.Range("H113").Value = 457.86667
$ws.Range("I113").Value = 394.0909
$ws.Range("J113").Value = 633.25
$ws.Range("K113").Value = 1182.2727
$ws.Range("L113").Value = 1899.75
$ws.Range("M113").Value = 987.7273
$ws.Range("N113").Value = -6239.75
